$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("D7").Value = "2016-27-12 06:27:03"
$wsZhCn.Range("E7").Value = "2016-03-12 06:27:00"
$wsDeDe.Range("E7").Value = "2016-03-12 06:27:03"
